# Added seasonality variable to data set
#
# The "Revenue" sheet had a "Winter North" Yes/No flag column (column E).
# It is replaced with a three-valued "Seasonality" column (High / Low /
# Medium) describing how busy each month is.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Revenue")

# Rename the column header.
$ws.Range("E1").Value = "Seasonality"

# New per-month seasonality values (rows 2-25 => 2022-01 .. 2023-12).
$seasonality = @(
    "High",   # 2022-01
    "High",   # 2022-02
    "High",   # 2022-03
    "Low",    # 2022-04
    "Low",    # 2022-05
    "Medium", # 2022-06
    "Medium", # 2022-07
    "Medium", # 2022-08
    "Low",    # 2022-09
    "Medium", # 2022-10
    "High",   # 2022-11
    "High",   # 2022-12
    "High",   # 2023-01
    "High",   # 2023-02
    "High",   # 2023-03
    "High",   # 2023-04
    "Low",    # 2023-05
    "Low",    # 2023-06
    "Medium", # 2023-07
    "Medium", # 2023-08
    "Low",    # 2023-09
    "Medium", # 2023-10
    "High",   # 2023-11
    "High"    # 2023-12
)

for ($i = 0; $i -lt $seasonality.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 5).Value = $seasonality[$i]
}

# Cosmetic: page setup + selected cell, matching the saved file state.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Activate()
$ws.Range("A2").Select() | Out-Null
